$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns to match the latest scrape.
# D-column price strings use "." as both thousands and decimal separators
# (as scraped from the source site), so each cell is forced to Text format
# before assignment and then restored to the default "Normal" style so that
# Excel does not silently reinterpret the text as a number.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.101.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.018.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "54.71"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.379"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.80%  "
$ws.Range("E10").Value = "  +2.07%  "
$ws.Range("E11").Value = "  -3.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.313.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.741"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.021.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.999.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.55%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0820"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.82%  "
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.125"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("E31").Value = "  -3.38%  "
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0615"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.03%  "
$ws.Range("E35").Value = "  -4.43%  "
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.22%  "
$ws.Range("E38").Value = "  -3.79%  "
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0218"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.477.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "16.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "95.29"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("E45").Value = "  -3.88%  "
$ws.Range("E46").Value = "  -3.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.51%  "
$ws.Range("E48").Value = "  -1.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.200.82"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.55%  "
$ws.Range("E51").Value = "  -10.44%  "
